# Applies the "error solve ifrs list" correction to the company_list sheet.
# Rows 2-6: financial figures recomputed/corrected (values shrink from raw KRW
#           thousands to a different, corrected scale); a handful of cells
#           (J, O, and some AD/AH cells) no longer have data and are cleared.
# Rows 7-9: only the id/name columns (A-C) remain; all metric columns (D:AJ)
#           are cleared out for these rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

    # Row 2
    $ws.Cells.Item(2, 4).Value = 5746
    $ws.Cells.Item(2, 5).Value = -26
    $ws.Cells.Item(2, 6).Value = -26
    $ws.Cells.Item(2, 7).Value = -41
    $ws.Cells.Item(2, 8).Value = -24
    $ws.Cells.Item(2, 9).Value = -24
    $ws.Cells.Item(2, 10).ClearContents()
    $ws.Cells.Item(2, 11).Value = 6292
    $ws.Cells.Item(2, 12).Value = 4669
    $ws.Cells.Item(2, 13).Value = 1623
    $ws.Cells.Item(2, 14).Value = 1623
    $ws.Cells.Item(2, 15).ClearContents()
    $ws.Cells.Item(2, 16).Value = 9
    $ws.Cells.Item(2, 17).Value = 130
    $ws.Cells.Item(2, 18).Value = -347
    $ws.Cells.Item(2, 19).Value = 165
    $ws.Cells.Item(2, 20).Value = 558
    $ws.Cells.Item(2, 21).Value = -429
    $ws.Cells.Item(2, 22).Value = 2790
    $ws.Cells.Item(2, 23).Value = -0.45
    $ws.Cells.Item(2, 24).Value = -0.42
    $ws.Cells.Item(2, 25).Value = -1.76
    $ws.Cells.Item(2, 26).Value = -0.43
    $ws.Cells.Item(2, 27).Value = 287.72
    $ws.Cells.Item(2, 28).Value = 20630.59
    $ws.Cells.Item(2, 29).Value = -180
    $ws.Cells.Item(2, 30).ClearContents()
    $ws.Cells.Item(2, 31).Value = 12165
    $ws.Cells.Item(2, 32).Value = 0
    $ws.Cells.Item(2, 33).Value = 0
    $ws.Cells.Item(2, 34).ClearContents()
    $ws.Cells.Item(2, 35).Value = 0
    $ws.Cells.Item(2, 36).Value = 13339817
    # Row 3
    $ws.Cells.Item(3, 4).Value = 6762
    $ws.Cells.Item(3, 5).Value = 454
    $ws.Cells.Item(3, 6).Value = 454
    $ws.Cells.Item(3, 7).Value = 198
    $ws.Cells.Item(3, 8).Value = 183
    $ws.Cells.Item(3, 9).Value = 183
    $ws.Cells.Item(3, 10).ClearContents()
    $ws.Cells.Item(3, 11).Value = 6555
    $ws.Cells.Item(3, 12).Value = 4665
    $ws.Cells.Item(3, 13).Value = 1890
    $ws.Cells.Item(3, 14).Value = 1890
    $ws.Cells.Item(3, 15).ClearContents()
    $ws.Cells.Item(3, 16).Value = 10
    $ws.Cells.Item(3, 17).Value = -131
    $ws.Cells.Item(3, 18).Value = -288
    $ws.Cells.Item(3, 19).Value = 532
    $ws.Cells.Item(3, 20).Value = 308
    $ws.Cells.Item(3, 21).Value = -439
    $ws.Cells.Item(3, 22).Value = 3363
    $ws.Cells.Item(3, 23).Value = 6.71
    $ws.Cells.Item(3, 24).Value = 2.71
    $ws.Cells.Item(3, 25).Value = 10.43
    $ws.Cells.Item(3, 26).Value = 2.85
    $ws.Cells.Item(3, 27).Value = 246.85
    $ws.Cells.Item(3, 28).Value = 18303.99
    $ws.Cells.Item(3, 29).Value = 1372
    $ws.Cells.Item(3, 30).ClearContents()
    $ws.Cells.Item(3, 31).Value = 12925
    $ws.Cells.Item(3, 32).Value = 0
    $ws.Cells.Item(3, 33).Value = 0
    $ws.Cells.Item(3, 34).ClearContents()
    $ws.Cells.Item(3, 35).Value = 0
    $ws.Cells.Item(3, 36).Value = 14621467
    # Row 4
    $ws.Cells.Item(4, 4).Value = 6771
    $ws.Cells.Item(4, 5).Value = 592
    $ws.Cells.Item(4, 6).Value = 592
    $ws.Cells.Item(4, 7).Value = 455
    $ws.Cells.Item(4, 8).Value = 347
    $ws.Cells.Item(4, 9).Value = 347
    $ws.Cells.Item(4, 10).ClearContents()
    $ws.Cells.Item(4, 11).Value = 6732
    $ws.Cells.Item(4, 12).Value = 3686
    $ws.Cells.Item(4, 13).Value = 3046
    $ws.Cells.Item(4, 14).Value = 3046
    $ws.Cells.Item(4, 15).ClearContents()
    $ws.Cells.Item(4, 16).Value = 108
    $ws.Cells.Item(4, 17).Value = 816
    $ws.Cells.Item(4, 18).Value = 147
    $ws.Cells.Item(4, 19).Value = -417
    $ws.Cells.Item(4, 20).Value = 125
    $ws.Cells.Item(4, 21).Value = 691
    $ws.Cells.Item(4, 22).Value = 2303
    $ws.Cells.Item(4, 23).Value = 8.75
    $ws.Cells.Item(4, 24).Value = 5.13
    $ws.Cells.Item(4, 25).Value = 14.06
    $ws.Cells.Item(4, 26).Value = 5.22
    $ws.Cells.Item(4, 27).Value = 121.02
    $ws.Cells.Item(4, 28).Value = 2713.08
    $ws.Cells.Item(4, 29).Value = 2106
    $ws.Cells.Item(4, 30).Value = 4.4
    $ws.Cells.Item(4, 31).Value = 13932
    $ws.Cells.Item(4, 32).Value = 0.66
    $ws.Cells.Item(4, 33).Value = 178
    $ws.Cells.Item(4, 34).Value = 1.92
    $ws.Cells.Item(4, 35).Value = 11.23
    $ws.Cells.Item(4, 36).Value = 21863268
    # Row 5
    $ws.Cells.Item(5, 4).Value = 6131
    $ws.Cells.Item(5, 5).Value = 319
    $ws.Cells.Item(5, 6).Value = 319
    $ws.Cells.Item(5, 7).Value = 312
    $ws.Cells.Item(5, 8).Value = 232
    $ws.Cells.Item(5, 9).Value = 234
    $ws.Cells.Item(5, 10).Value = -1
    $ws.Cells.Item(5, 11).Value = 6756
    $ws.Cells.Item(5, 12).Value = 3519
    $ws.Cells.Item(5, 13).Value = 3237
    $ws.Cells.Item(5, 14).Value = 3232
    $ws.Cells.Item(5, 15).Value = 5
    $ws.Cells.Item(5, 16).Value = 109
    $ws.Cells.Item(5, 17).Value = 283
    $ws.Cells.Item(5, 18).Value = -132
    $ws.Cells.Item(5, 19).Value = 81
    $ws.Cells.Item(5, 20).Value = 196
    $ws.Cells.Item(5, 21).Value = 87
    $ws.Cells.Item(5, 22).Value = 2351
    $ws.Cells.Item(5, 23).Value = 5.2
    $ws.Cells.Item(5, 24).Value = 3.79
    $ws.Cells.Item(5, 25).Value = 7.45
    $ws.Cells.Item(5, 26).Value = 3.45
    $ws.Cells.Item(5, 27).Value = 108.74
    $ws.Cells.Item(5, 28).Value = 2875.56
    $ws.Cells.Item(5, 29).Value = 1069
    $ws.Cells.Item(5, 30).Value = 8.460000000000001
    $ws.Cells.Item(5, 31).Value = 14783
    $ws.Cells.Item(5, 32).Value = 0.61
    $ws.Cells.Item(5, 33).Value = 100
    $ws.Cells.Item(5, 34).Value = 1.1
    $ws.Cells.Item(5, 35).Value = 9.35
    $ws.Cells.Item(5, 36).Value = 21863268
    # Row 6
    $ws.Cells.Item(6, 4).Value = 6009
    $ws.Cells.Item(6, 5).Value = 105
    $ws.Cells.Item(6, 6).Value = 105
    $ws.Cells.Item(6, 7).Value = 36
    $ws.Cells.Item(6, 8).Value = 41
    $ws.Cells.Item(6, 9).Value = 42
    $ws.Cells.Item(6, 11).Value = 7053
    $ws.Cells.Item(6, 12).Value = 3847
    $ws.Cells.Item(6, 13).Value = 3206
    $ws.Cells.Item(6, 14).Value = 3203
    $ws.Cells.Item(6, 16).Value = 109
    $ws.Cells.Item(6, 17).Value = 195
    $ws.Cells.Item(6, 18).Value = -1012
    $ws.Cells.Item(6, 19).Value = 184
    $ws.Cells.Item(6, 20).Value = 975
    $ws.Cells.Item(6, 21).Value = -780
    $ws.Cells.Item(6, 22).Value = 2589
    $ws.Cells.Item(6, 23).Value = 1.75
    $ws.Cells.Item(6, 24).Value = 0.68
    $ws.Cells.Item(6, 25).Value = 1.31
    $ws.Cells.Item(6, 26).Value = 0.59
    $ws.Cells.Item(6, 27).Value = 119.98
    $ws.Cells.Item(6, 28).Value = 2903.55
    $ws.Cells.Item(6, 29).Value = 192
    $ws.Cells.Item(6, 30).Value = 32.33
    $ws.Cells.Item(6, 31).Value = 14909
    $ws.Cells.Item(6, 32).Value = 0.42
    $ws.Cells.Item(6, 33).Value = 100
    $ws.Cells.Item(6, 34).Value = 1.61
    $ws.Cells.Item(6, 35).Value = 23.29
    $ws.Cells.Item(6, 36).Value = 21863268


    # Rows 7-9: clear every metric cell (D:AJ), keep A/B/C (id, ticker, name)
    $ws.Range("D7:AJ9").ClearContents()
